$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{Row=1;  Col=1; Text="87÷5=17, 2"},
    @{Row=1;  Col=2; Text="23÷7=3, 2"},
    @{Row=1;  Col=3; Text="43÷4=10, 3"},
    @{Row=1;  Col=4; Text="83÷7=11, 6"},
    @{Row=1;  Col=5; Text="51÷9=5, 6"},

    @{Row=5;  Col=1; Text="36÷4=9, 0"},
    @{Row=5;  Col=2; Text="27÷9=3, 0"},
    @{Row=5;  Col=3; Text="34÷9=3, 7"},
    @{Row=5;  Col=4; Text="48÷2=24, 0"},
    @{Row=5;  Col=5; Text="89÷6=14, 5"},

    @{Row=9;  Col=1; Text="91÷6=15, 1"},
    @{Row=9;  Col=2; Text="52÷6=8, 4"},
    @{Row=9;  Col=3; Text="60÷2=30, 0"},
    @{Row=9;  Col=4; Text="74÷7=10, 4"},
    @{Row=9;  Col=5; Text="47÷9=5, 2"},

    @{Row=13; Col=1; Text="73÷5=14, 3"},
    @{Row=13; Col=2; Text="27÷6=4, 3"},
    @{Row=13; Col=3; Text="60÷2=30, 0"},
    @{Row=13; Col=4; Text="18÷3=6, 0"},
    @{Row=13; Col=5; Text="84÷5=16, 4"},

    @{Row=17; Col=1; Text="12÷5=2, 2"},
    @{Row=17; Col=2; Text="58÷5=11, 3"},
    @{Row=17; Col=3; Text="83÷8=10, 3"},
    @{Row=17; Col=4; Text="64÷3=21, 1"},
    @{Row=17; Col=5; Text="14÷3=4, 2"}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
